$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 4th execution case (old row 13: 3/Pdt.Eks/2022/PA.Tte dated 07-08/11/2022 ... 16/11/2022)
# is removed as a distinct row; its data is folded up into row 12, row 12's old data moves to
# row 11, and row 11's old data moves to row 10 (filling the previously-empty G:I columns).
# Deleting row 13 outright shifts the footer (signature block) rows up by one to match.
$ws.Rows(13).Delete()

# Row 10 gains the "Penetapan/Pelaksanaan Sita Eksekusi" + "Penetapan Eksekusi" dates that
# used to live on row 11 (G:I) of the old layout.
$ws.Range("G10").Value = "31/10/2022"
$ws.Range("H10").Value = "11/11/2022"
$ws.Range("I10").Value = "30/12/2022"

# Row 11 now carries what used to be row 12's case, with G:I cleared and a new "Bergantung"
# date (J11) added.
$ws.Range("B11").Value = "2/Pdt.Eks/2021/PA.Tte"
$ws.Range("C11").Value = "446/Pdt.G/2020/PA.Tte"
$ws.Range("D11").Value = "16/06/2021"
$ws.Range("E11").Value = "18/06/2021"
$ws.Range("F11").Value = "28/06/2021"
$ws.Range("G11").ClearContents()
$ws.Range("H11").ClearContents()
$ws.Range("I11").ClearContents()
$ws.Range("J11").Value = "11/01/2023"

# Row 12 now carries what used to be row 13's case.
$ws.Range("B12").Value = "3/Pdt.Eks/2022/PA.Tte"
$ws.Range("C12").ClearContents()
$ws.Range("D12").Value = "07/11/2022"
$ws.Range("E12").Value = "08/11/2022"
$ws.Range("F12").Value = "16/11/2022"

# Signature block date updated (row 16 -> row 15 after the row delete).
$ws.Range("J15").Value = "Ternate , 04 September 2023"

# Keep the saved selection in sync with the new layout.
$ws.Range("C16").Select()
